# Weekly update: insert a new price record for Ciboulette (Femacal de La
# Calera) as the most recent observation. The sheet keeps its rows ordered
# with new weekly entries inserted at row 17, pushing the older history
# (previously rows 17-200) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 - this shifts rows 17:200 down to 18:201
# and carries the dimension/date-column formatting along with it.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly observation. All
# fields mirror the prior row 17 record except the date (column D), which
# is the new, more recent sample date.
$ws.Range("A17").Value = 3
$ws.Range("B17").Value = "Femacal de La Calera"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = 44503
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 100112039
$ws.Range("G17").Value = "Ciboulette"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 160
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = 1500
$ws.Range("N17").Value = "$/docena de atados"
$ws.Range("O17").Value = "Provincia de Quillota"
$ws.Range("P17").Value = 500
$ws.Range("Q17").Value = 3
$ws.Range("R17").Value = "Hortaliza"
